$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Dopamine Detox: A Short Guide to Remove Distractions and Get Your Brain to Do Hard Things"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "289."
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "4.4 out of 5 stars"
$ws.Range("D2").Value = "Thibaut Meurisse"

# Row 3
$ws.Range("A3").Value = "Don't Believe Everything You Think (English)"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "184."
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "4.5 out of 5 stars"
$ws.Range("D3").Value = "Joseph Nguyen"

# Row 4
$ws.Range("A4").Value = "The Psychology of Money"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "160."
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "4.6 out of 5 stars"
$ws.Range("D4").Value = "Morgan Housel"

# Row 5 (new)
$ws.Range("A5").Value = "White Nights – Fyodor Dostoyevsky | A Million-Copy Bestseller | A Timeless Story of Love, Longing & Solitude – Penguin Classics"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "89."
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "4.6 out of 5 stars"
$ws.Range("D5").Value = "Fyodor Dostoyevsky"

# Row 6 (new)
$ws.Range("A6").Value = "The Art of Being Alone: Loneliness Was My Cage, Solitude Is My Home (English)"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "199."
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "4.4 out of 5 stars"
$ws.Range("D6").Value = "Renuka Gavrani"
